$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Result")
$ws.Activate()

# A new scenario ("View Guardian") was added to the checklist, so insert a
# fresh row above the old row 4 ("Manage Patient" / "Add Patient"); everything
# below shifts down by one and the header formula's range auto-extends
# (C993 -> C994).
$ws.Rows.Item(4).Insert()

# Give the new row the same banded border used by every other data row (it
# is not inherited automatically on an inserted blank row).
foreach ($col in @("A", "B")) {
    $edge = $ws.Range($col + "4")
    $edge.Borders.Item(7).LineStyle = 1
    $edge.Borders.Item(7).Weight = -4138
    $edge.Borders.Item(7).Color = 13421772
    $edge.Borders.Item(8).LineStyle = 1
    $edge.Borders.Item(8).Weight = -4138
    $edge.Borders.Item(8).Color = 13421772
    $edge.Borders.Item(9).LineStyle = 1
    $edge.Borders.Item(9).Weight = -4138
    $edge.Borders.Item(9).Color = 13421772
    $edge.Borders.Item(10).LineStyle = 1
    $edge.Borders.Item(10).Weight = -4138
    $edge.Borders.Item(10).Color = 13421772
}
$ws.Range("B4").WrapText = $true
$ws.Range("D4").WrapText = $true

# New scenario name + not-yet-validated result.
$ws.Range("B4").Value = "View Guardian"
$ws.Range("C4").Value = $false

# This test run: none of the scenarios passed validation, so every result
# checkbox in the column (including the new row) is reset to FALSE.
$ws.Range("C2:C12").Value = $false

# Leave the selection where the author left it.
[void]$ws.Range("D10").Select()

Write-Output "Added 'View Guardian' test result row and reset validation results"
